$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(53, 8).Value = 319.875  # H53: 448.5 -> 319.875
$ws.Cells.Item(53, 9).Value = 312.8  # I53: 498.5 -> 312.8
$ws.Cells.Item(53, 10).Value = 331.66666  # J53: 398.5 -> 331.66666
$ws.Cells.Item(53, 11).Value = 312.8  # K53: 498.5 -> 312.8
$ws.Cells.Item(53, 12).Value = 331.66666  # L53: 398.5 -> 331.66666
$ws.Cells.Item(53, 13).Value = 324.2  # M53: 138.5 -> 324.2
$ws.Cells.Item(53, 14).Value = -1605.66666  # N53: -1672.5 -> -1605.66666
$ws.Cells.Item(62, 8).Value = 1666979.1  # H62: 3333333.2 -> 1666979.1
$ws.Cells.Item(62, 9).Value = 2000241.8  # I62: 3333333.2 -> 2000241.8
$ws.Cells.Item(62, 10).Value = 666  # J62: 0 -> 666
$ws.Cells.Item(62, 11).Value = 2000241.8  # K62: 3333333.2 -> 2000241.8
$ws.Cells.Item(62, 12).Value = 666  # L62: 0 -> 666
$ws.Cells.Item(62, 13).Value = -1999617.8  # M62: -3332709.2 -> -1999617.8
$ws.Cells.Item(62, 14).Value = -1914  # N62: None -> -1914
$ws.Cells.Item(65, 8).Value = 1666979.1  # H65: 3333333.2 -> 1666979.1
$ws.Cells.Item(65, 9).Value = 2000241.8  # I65: 3333333.2 -> 2000241.8
$ws.Cells.Item(65, 10).Value = 666  # J65: 0 -> 666
$ws.Cells.Item(65, 11).Value = 10001209  # K65: 16666666 -> 10001209
$ws.Cells.Item(65, 12).Value = 3330  # L65: 0 -> 3330
$ws.Cells.Item(65, 13).Value = -9998089  # M65: -16663546 -> -9998089
$ws.Cells.Item(65, 14).Value = -9570  # N65: None -> -9570
$ws.Cells.Item(86, 8).Value = 151008.3  # H86: 145631.22 -> 151008.3
$ws.Cells.Item(86, 9).Value = 288693.44  # I86: 269477.2 -> 288693.44
$ws.Cells.Item(86, 11).Value = 288693.44  # K86: 269477.2 -> 288693.44
$ws.Cells.Item(86, 13).Value = -287570.44  # M86: -268354.2 -> -287570.44
$ws.Cells.Item(89, 8).Value = 151008.3  # H89: 145631.22 -> 151008.3
$ws.Cells.Item(89, 9).Value = 288693.44  # I89: 269477.2 -> 288693.44
$ws.Cells.Item(89, 11).Value = 1443467.2  # K89: 1347386 -> 1443467.2
$ws.Cells.Item(89, 13).Value = -1437851.2  # M89: -1341770 -> -1437851.2
$ws.Cells.Item(99, 8).Value = 6499  # H99: 1538.5 -> 6499
$ws.Cells.Item(99, 9).Value = 0  # I99: 160 -> 0
$ws.Cells.Item(99, 10).Value = 6499  # J99: 1998 -> 6499
$ws.Cells.Item(99, 11).Value = 0  # K99: 480 -> 0
$ws.Cells.Item(99, 12).Value = 19497  # L99: 5994 -> 19497
$ws.Cells.Item(99, 13).Value = ""  # M99: clear (was 1018)
$ws.Cells.Item(99, 14).Value = -22493  # N99: -8990 -> -22493
$ws.Cells.Item(100, 8).Value = 2416.5  # H100: 2625 -> 2416.5
$ws.Cells.Item(100, 9).Value = 2047.5714  # I100: 2229.3333 -> 2047.5714
$ws.Cells.Item(100, 11).Value = 2047.5714  # K100: 2229.3333 -> 2047.5714
$ws.Cells.Item(100, 13).Value = -1506.5714  # M100: -1688.3333 -> -1506.5714
$ws.Cells.Item(106, 8).Value = 1818.875  # H106: 5102.7334 -> 1818.875
$ws.Cells.Item(106, 9).Value = 1818.875  # I106: 1837.8889 -> 1818.875
$ws.Cells.Item(106, 10).Value = 0  # J106: 10000 -> 0
$ws.Cells.Item(106, 11).Value = 1818.875  # K106: 1837.8889 -> 1818.875
$ws.Cells.Item(106, 12).Value = 0  # L106: 10000 -> 0
$ws.Cells.Item(106, 13).Value = -1187.875  # M106: -1206.8889 -> -1187.875
$ws.Cells.Item(106, 14).Value = ""  # N106: clear (was -11262)
$ws.Cells.Item(112, 8).Value = 3629.3704  # H112: 3499 -> 3629.3704
$ws.Cells.Item(112, 10).Value = 4028  # J112: 3979.389 -> 4028
$ws.Cells.Item(112, 12).Value = 12084  # L112: 11938.167 -> 12084
$ws.Cells.Item(112, 14).Value = -14300  # N112: -14154.167 -> -14300

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2806.9048  # H2: 2578.9167 -> 2806.9048
$ws.Cells.Item(2, 9).Value = 2016.3334  # I2: 2251.5386 -> 2016.3334
$ws.Cells.Item(2, 10).Value = 4783.3335  # J2: 2965.818 -> 4783.3335
$ws.Cells.Item(2, 11).Value = 2016.3334  # K2: 2251.5386 -> 2016.3334
$ws.Cells.Item(2, 12).Value = 4783.3335  # L2: 2965.818 -> 4783.3335
$ws.Cells.Item(2, 13).Value = -1903.3334  # M2: -2138.5386 -> -1903.3334
$ws.Cells.Item(2, 14).Value = -5009.3335  # N2: -3191.818 -> -5009.3335
$ws.Cells.Item(74, 8).Value = 1983  # H74: 1931.2941 -> 1983
$ws.Cells.Item(74, 9).Value = 1303.0714  # I74: 1265.8966 -> 1303.0714
$ws.Cells.Item(74, 11).Value = 1303.0714  # K74: 1265.8966 -> 1303.0714
$ws.Cells.Item(74, 13).Value = -429.0714  # M74: -391.8966 -> -429.0714
$ws.Cells.Item(77, 8).Value = 1983  # H77: 1931.2941 -> 1983
$ws.Cells.Item(77, 9).Value = 1303.0714  # I77: 1265.8966 -> 1303.0714
$ws.Cells.Item(77, 11).Value = 6515.357  # K77: 6329.483 -> 6515.357
$ws.Cells.Item(77, 13).Value = -2147.357  # M77: -1961.483 -> -2147.357
$ws.Cells.Item(86, 8).Value = 11475  # H86: 11316.667 -> 11475
$ws.Cells.Item(86, 10).Value = 11475  # J86: 11316.667 -> 11475
$ws.Cells.Item(86, 12).Value = 11475  # L86: 11316.667 -> 11475
$ws.Cells.Item(86, 14).Value = -13847  # N86: -13688.667 -> -13847
$ws.Cells.Item(89, 8).Value = 11475  # H89: 11316.667 -> 11475
$ws.Cells.Item(89, 10).Value = 11475  # J89: 11316.667 -> 11475
$ws.Cells.Item(89, 12).Value = 34425  # L89: 33950.001 -> 34425
$ws.Cells.Item(89, 14).Value = -46281  # N89: -45806.001 -> -46281
$ws.Cells.Item(116, 8).Value = 2806.9048  # H116: 2578.9167 -> 2806.9048
$ws.Cells.Item(116, 9).Value = 2016.3334  # I116: 2251.5386 -> 2016.3334
$ws.Cells.Item(116, 10).Value = 4783.3335  # J116: 2965.818 -> 4783.3335
$ws.Cells.Item(116, 11).Value = 2016.3334  # K116: 2251.5386 -> 2016.3334
$ws.Cells.Item(116, 12).Value = 4783.3335  # L116: 2965.818 -> 4783.3335
$ws.Cells.Item(116, 13).Value = 277.6666  # M116: 42.46140000000014 -> 277.6666
$ws.Cells.Item(116, 14).Value = -9371.333500000001  # N116: -7553.818 -> -9371.333500000001
$ws.Cells.Item(132, 8).Value = 3175.9062  # H132: 3278.8125 -> 3175.9062
$ws.Cells.Item(132, 9).Value = 3184.3215  # I132: 3301.9285 -> 3184.3215
$ws.Cells.Item(132, 11).Value = 9552.9645  # K132: 9905.7855 -> 9552.9645
$ws.Cells.Item(132, 13).Value = -7022.9645  # M132: -7375.7855 -> -7022.9645
$ws.Cells.Item(133, 8).Value = 51323.6  # H133: 44470.668 -> 51323.6
$ws.Cells.Item(133, 10).Value = 48111  # J133: 39999.875 -> 48111
$ws.Cells.Item(133, 12).Value = 48111  # L133: 39999.875 -> 48111
$ws.Cells.Item(133, 14).Value = -53171  # N133: -45059.875 -> -53171

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(16, 8).Value = 10000  # H16: 908 -> 10000
$ws.Cells.Item(16, 9).Value = 10000  # I16: 908 -> 10000
$ws.Cells.Item(16, 11).Value = 10000  # K16: 908 -> 10000
$ws.Cells.Item(16, 13).Value = -9830  # M16: -738 -> -9830
$ws.Cells.Item(86, 8).Value = 2691.4583  # H86: 2751.84 -> 2691.4583
$ws.Cells.Item(86, 9).Value = 2290.1667  # I86: 2437.1538 -> 2290.1667
$ws.Cells.Item(86, 11).Value = 2290.1667  # K86: 2437.1538 -> 2290.1667
$ws.Cells.Item(86, 13).Value = -1167.1667  # M86: -1314.1538 -> -1167.1667
$ws.Cells.Item(89, 8).Value = 2691.4583  # H89: 2751.84 -> 2691.4583
$ws.Cells.Item(89, 9).Value = 2290.1667  # I89: 2437.1538 -> 2290.1667
$ws.Cells.Item(89, 11).Value = 11450.8335  # K89: 12185.769 -> 11450.8335
$ws.Cells.Item(89, 13).Value = -5834.833500000001  # M89: -6569.769 -> -5834.833500000001
$ws.Cells.Item(105, 8).Value = 2500.3635  # H105: 2600.4 -> 2500.3635
$ws.Cells.Item(105, 9).Value = 2500.3635  # I105: 2600.4 -> 2500.3635
$ws.Cells.Item(105, 11).Value = 2500.3635  # K105: 2600.4 -> 2500.3635
$ws.Cells.Item(105, 13).Value = -753.3634999999999  # M105: -853.4000000000001 -> -753.3634999999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(5, 8).Value = 12459.4  # H5: 24999.5 -> 12459.4
$ws.Cells.Item(5, 9).Value = 750  # I5: 20000 -> 750
$ws.Cells.Item(5, 10).Value = 20265.666  # J5: 29999 -> 20265.666
$ws.Cells.Item(5, 11).Value = 750  # K5: 20000 -> 750
$ws.Cells.Item(5, 12).Value = 20265.666  # L5: 29999 -> 20265.666
$ws.Cells.Item(5, 13).Value = -638  # M5: -19888 -> -638
$ws.Cells.Item(5, 14).Value = -20489.666  # N5: -30223 -> -20489.666
$ws.Cells.Item(8, 8).Value = 29999.5  # H8: 99999 -> 29999.5
$ws.Cells.Item(8, 10).Value = 29999.5  # J8: 99999 -> 29999.5
$ws.Cells.Item(8, 12).Value = 29999.5  # L8: 99999 -> 29999.5
$ws.Cells.Item(8, 14).Value = -30279.5  # N8: -100279 -> -30279.5
$ws.Cells.Item(15, 8).Value = 144  # H15: 646 -> 144
$ws.Cells.Item(15, 9).Value = 144  # I15: 188 -> 144
$ws.Cells.Item(15, 10).Value = 0  # J15: 875 -> 0
$ws.Cells.Item(15, 11).Value = 144  # K15: 188 -> 144
$ws.Cells.Item(15, 12).Value = 0  # L15: 875 -> 0
$ws.Cells.Item(15, 13).Value = 26  # M15: -18 -> 26
$ws.Cells.Item(15, 14).Value = ""  # N15: clear (was -1215)
$ws.Cells.Item(41, 8).Value = 28916.666  # H41: 750 -> 28916.666
$ws.Cells.Item(41, 10).Value = 43000  # J41: 0 -> 43000
$ws.Cells.Item(41, 12).Value = 43000  # L41: 0 -> 43000
$ws.Cells.Item(41, 14).Value = -43856  # N41: None -> -43856
$ws.Cells.Item(50, 8).Value = 21500  # H50: 22166.666 -> 21500
$ws.Cells.Item(50, 10).Value = 22000  # J50: 23250 -> 22000
$ws.Cells.Item(50, 12).Value = 22000  # L50: 23250 -> 22000
$ws.Cells.Item(50, 14).Value = -23250  # N50: -24500 -> -23250
$ws.Cells.Item(51, 8).Value = 45000  # H51: 0 -> 45000
$ws.Cells.Item(51, 10).Value = 45000  # J51: 0 -> 45000
$ws.Cells.Item(51, 12).Value = 45000  # L51: 0 -> 45000
$ws.Cells.Item(51, 14).Value = -46472  # N51: None -> -46472
$ws.Cells.Item(59, 8).Value = 106856.43  # H59: 108332.5 -> 106856.43
$ws.Cells.Item(59, 10).Value = 106856.43  # J59: 108332.5 -> 106856.43
$ws.Cells.Item(59, 12).Value = 106856.43  # L59: 108332.5 -> 106856.43
$ws.Cells.Item(59, 14).Value = -109146.43  # N59: -110622.5 -> -109146.43
$ws.Cells.Item(60, 8).Value = 37989.75  # H60: 37999.5 -> 37989.75
$ws.Cells.Item(60, 10).Value = 37986.332  # J60: 37999 -> 37986.332
$ws.Cells.Item(60, 12).Value = 37986.332  # L60: 37999 -> 37986.332
$ws.Cells.Item(60, 14).Value = -39008.332  # N60: -39021 -> -39008.332
$ws.Cells.Item(61, 8).Value = 45000  # H61: 0 -> 45000
$ws.Cells.Item(61, 10).Value = 45000  # J61: 0 -> 45000
$ws.Cells.Item(61, 12).Value = 45000  # L61: 0 -> 45000
$ws.Cells.Item(61, 14).Value = -45696  # N61: None -> -45696
$ws.Cells.Item(62, 8).Value = 15087.8  # H62: 8617.375 -> 15087.8
$ws.Cells.Item(62, 9).Value = 17999  # I62: 3166.3333 -> 17999
$ws.Cells.Item(62, 10).Value = 14360  # J62: 11888 -> 14360
$ws.Cells.Item(62, 11).Value = 17999  # K62: 3166.3333 -> 17999
$ws.Cells.Item(62, 12).Value = 14360  # L62: 11888 -> 14360
$ws.Cells.Item(62, 13).Value = -17375  # M62: -2542.3333 -> -17375
$ws.Cells.Item(62, 14).Value = -15608  # N62: -13136 -> -15608
$ws.Cells.Item(65, 8).Value = 15087.8  # H65: 8617.375 -> 15087.8
$ws.Cells.Item(65, 9).Value = 17999  # I65: 3166.3333 -> 17999
$ws.Cells.Item(65, 10).Value = 14360  # J65: 11888 -> 14360
$ws.Cells.Item(65, 11).Value = 89995  # K65: 15831.6665 -> 89995
$ws.Cells.Item(65, 12).Value = 71800  # L65: 59440 -> 71800
$ws.Cells.Item(65, 13).Value = -86875  # M65: -12711.6665 -> -86875
$ws.Cells.Item(65, 14).Value = -78040  # N65: -65680 -> -78040
$ws.Cells.Item(74, 8).Value = 30314  # H74: 0 -> 30314
$ws.Cells.Item(74, 10).Value = 30314  # J74: 0 -> 30314
$ws.Cells.Item(74, 12).Value = 30314  # L74: 0 -> 30314
$ws.Cells.Item(74, 14).Value = -32062  # N74: None -> -32062
$ws.Cells.Item(77, 8).Value = 30314  # H77: 0 -> 30314
$ws.Cells.Item(77, 10).Value = 30314  # J77: 0 -> 30314
$ws.Cells.Item(77, 12).Value = 90942  # L77: 0 -> 90942
$ws.Cells.Item(77, 14).Value = -99678  # N77: None -> -99678
$ws.Cells.Item(107, 8).Value = 1180.4286  # H107: 757.2381 -> 1180.4286
$ws.Cells.Item(107, 9).Value = 1189.5  # I107: 745.15 -> 1189.5
$ws.Cells.Item(107, 11).Value = 1189.5  # K107: 745.15 -> 1189.5
$ws.Cells.Item(107, 13).Value = 730.5  # M107: 1174.85 -> 730.5
$ws.Cells.Item(134, 8).Value = 3830.2258  # H134: 3911.4666 -> 3830.2258
$ws.Cells.Item(134, 9).Value = 3824.9312  # I134: 3911.7856 -> 3824.9312
$ws.Cells.Item(134, 11).Value = 11474.7936  # K134: 11735.3568 -> 11474.7936
$ws.Cells.Item(134, 13).Value = -8939.793600000001  # M134: -9200.356800000001 -> -8939.793600000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(98, 8).Value = 1432.909  # H98: 1489.5834 -> 1432.909
$ws.Cells.Item(98, 9).Value = 1183  # I98: 888 -> 1183
$ws.Cells.Item(98, 10).Value = 1526.625  # J98: 1790.375 -> 1526.625
$ws.Cells.Item(98, 11).Value = 3549  # K98: 2664 -> 3549
$ws.Cells.Item(98, 12).Value = 4579.875  # L98: 5371.125 -> 4579.875
$ws.Cells.Item(98, 13).Value = -2051  # M98: -1166 -> -2051
$ws.Cells.Item(98, 14).Value = -7575.875  # N98: -8367.125 -> -7575.875

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(103, 8).Value = 49299  # H103: 0 -> 49299
$ws.Cells.Item(103, 10).Value = 49299  # J103: 0 -> 49299
$ws.Cells.Item(103, 12).Value = 49299  # L103: 0 -> 49299
$ws.Cells.Item(103, 14).Value = -51643  # N103: None -> -51643
$ws.Cells.Item(105, 8).Value = 46164.332  # H105: 46832.5 -> 46164.332
$ws.Cells.Item(105, 10).Value = 46164.332  # J105: 46832.5 -> 46164.332
$ws.Cells.Item(105, 12).Value = 46164.332  # L105: 46832.5 -> 46164.332
$ws.Cells.Item(105, 14).Value = -53152.332  # N105: -53820.5 -> -53152.332

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 9000  # H68: 4596.6665 -> 9000
$ws.Cells.Item(68, 9).Value = 9000  # I68: 4596.6665 -> 9000
$ws.Cells.Item(68, 11).Value = 9000  # K68: 4596.6665 -> 9000
$ws.Cells.Item(68, 13).Value = -8251  # M68: -3847.6665 -> -8251
$ws.Cells.Item(71, 8).Value = 9000  # H71: 4596.6665 -> 9000
$ws.Cells.Item(71, 9).Value = 9000  # I71: 4596.6665 -> 9000
$ws.Cells.Item(71, 11).Value = 45000  # K71: 22983.3325 -> 45000
$ws.Cells.Item(71, 13).Value = -41256  # M71: -19239.3325 -> -41256
$ws.Cells.Item(101, 8).Value = 24178.6  # H101: 24299.334 -> 24178.6
$ws.Cells.Item(101, 10).Value = 24178.6  # J101: 24299.334 -> 24178.6
$ws.Cells.Item(101, 12).Value = 24178.6  # L101: 24299.334 -> 24178.6
$ws.Cells.Item(101, 14).Value = -30668.6  # N101: -30789.334 -> -30668.6
$ws.Cells.Item(140, 8).Value = 77139.45  # H140: 80610.5 -> 77139.45
$ws.Cells.Item(140, 10).Value = 71358.39999999999  # J140: 74572.78 -> 71358.39999999999
$ws.Cells.Item(140, 12).Value = 71358.39999999999  # L140: 74572.78 -> 71358.39999999999
$ws.Cells.Item(140, 14).Value = -81718.39999999999  # N140: -84932.78 -> -81718.39999999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 4928.136  # H122: 4948.5454 -> 4928.136
$ws.Cells.Item(122, 9).Value = 4231.722  # I122: 4374.7646 -> 4231.722
$ws.Cells.Item(122, 10).Value = 8062  # J122: 6899.4 -> 8062
$ws.Cells.Item(122, 11).Value = 12695.166  # K122: 13124.2938 -> 12695.166
$ws.Cells.Item(122, 12).Value = 24186  # L122: 20698.2 -> 24186
$ws.Cells.Item(122, 13).Value = -10245.166  # M122: -10674.2938 -> -10245.166
$ws.Cells.Item(122, 14).Value = -29086  # N122: -25598.2 -> -29086
